{"js": "// Replace each three-digit-by-one-digit multiplication equation in the\n// document's table cells with its updated counterpart, per the commit diff.\n// Every \"old\" equation string below is unique in the document, so an exact\n// (case-sensitive, non-wildcard) search-and-replace is unambiguous.\nconst replacements = [\n  [\"464\u00d76=2784\", \"289\u00d73=867\"],\n  [\"619\u00d79=5571\", \"725\u00d78=5800\"],\n  [\"913\u00d75=4565\", \"612\u00d73=1836\"],\n  [\"555\u00d73=1665\", \"957\u00d76=5742\"],\n  [\"589\u00d73=1767\", \"528\u00d76=3168\"],\n  [\"243\u00d74=972\", \"257\u00d72=514\"],\n  [\"280\u00d77=1960\", \"805\u00d79=7245\"],\n  [\"281\u00d76=1686\", \"519\u00d74=2076\"],\n  [\"288\u00d74=1152\", \"783\u00d75=3915\"],\n  [\"853\u00d79=7677\", \"378\u00d77=2646\"],\n  [\"391\u00d75=1955\", \"597\u00d78=4776\"],\n  [\"886\u00d77=6202\", \"196\u00d74=784\"],\n  [\"194\u00d77=1358\", \"635\u00d79=5715\"],\n  [\"817\u00d75=4085\", \"896\u00d75=4480\"],\n  [\"353\u00d77=2471\", \"948\u00d75=4740\"],\n  [\"186\u00d76=1116\", \"435\u00d76=2610\"],\n  [\"558\u00d79=5022\", \"941\u00d75=4705\"],\n  [\"706\u00d76=4236\", \"185\u00d78=1480\"],\n  [\"179\u00d72=358\", \"553\u00d77=3871\"],\n  [\"485\u00d76=2910\", \"903\u00d76=5418\"],\n  [\"331\u00d78=2648\", \"105\u00d73=315\"],\n  [\"611\u00d73=1833\", \"831\u00d72=1662\"],\n  [\"168\u00d77=1176\", \"199\u00d74=796\"],\n  [\"386\u00d77=2702\", \"849\u00d74=3396\"],\n  [\"568\u00d76=3408\", \"795\u00d75=3975\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication equation in the\n# document's table cells with its updated counterpart, per the commit diff.\n# Every \"old\" equation string is unique in the document, so an exact\n# (case-sensitive, non-wildcard) Find/Replace is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"464\u00d76=2784\", \"289\u00d73=867\"),\n    @(\"619\u00d79=5571\", \"725\u00d78=5800\"),\n    @(\"913\u00d75=4565\", \"612\u00d73=1836\"),\n    @(\"555\u00d73=1665\", \"957\u00d76=5742\"),\n    @(\"589\u00d73=1767\", \"528\u00d76=3168\"),\n    @(\"243\u00d74=972\", \"257\u00d72=514\"),\n    @(\"280\u00d77=1960\", \"805\u00d79=7245\"),\n    @(\"281\u00d76=1686\", \"519\u00d74=2076\"),\n    @(\"288\u00d74=1152\", \"783\u00d75=3915\"),\n    @(\"853\u00d79=7677\", \"378\u00d77=2646\"),\n    @(\"391\u00d75=1955\", \"597\u00d78=4776\"),\n    @(\"886\u00d77=6202\", \"196\u00d74=784\"),\n    @(\"194\u00d77=1358\", \"635\u00d79=5715\"),\n    @(\"817\u00d75=4085\", \"896\u00d75=4480\"),\n    @(\"353\u00d77=2471\", \"948\u00d75=4740\"),\n    @(\"186\u00d76=1116\", \"435\u00d76=2610\"),\n    @(\"558\u00d79=5022\", \"941\u00d75=4705\"),\n    @(\"706\u00d76=4236\", \"185\u00d78=1480\"),\n    @(\"179\u00d72=358\", \"553\u00d77=3871\"),\n    @(\"485\u00d76=2910\", \"903\u00d76=5418\"),\n    @(\"331\u00d78=2648\", \"105\u00d73=315\"),\n    @(\"611\u00d73=1833\", \"831\u00d72=1662\"),\n    @(\"168\u00d77=1176\", \"199\u00d74=796\"),\n    @(\"386\u00d77=2702\", \"849\u00d74=3396\"),\n    @(\"568\u00d76=3408\", \"795\u00d75=3975\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
